$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "S" note in column H (row 22). Doing this before writing
# the new row-31 topic lets the underlying shared-string slot be reclaimed
# and reused in place, instead of appending a duplicate string entry.
$ws.Range("H22").ClearContents()

# Fix a typo'd date (was year 2011, should be 2019-10-25 like its neighbours).
$ws.Range("A30").Value = 43763

# Add the new "Front page + requirements + document structure" entry as the
# final row of the third effort table.
$ws.Range("A31").Value = 43767
$ws.Range("B31").Value = "Front page + requirements + document structure"
$ws.Range("C31").Value = 3
$ws.Rows.Item(31).RowHeight = 43.5

# C32 holds =SUM(C26:C31); it recalculates automatically to include the
# newly-populated C31 (11 -> 14).

# Update the view's selection to match where the author left the cursor.
$ws.Range("E17").Select() | Out-Null
